$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the date-cell number format used by column D before inserting rows
$dateFmt = $ws.Cells.Item(186, 4).NumberFormat

# Insert 3 new rows above the current row 186 (old data 186-207 shifts down to 189-210)
$ws.Rows.Item(186).Insert()
$ws.Rows.Item(186).Insert()
$ws.Rows.Item(186).Insert()

# New row 186: Espárragos, Banquete
$ws.Cells.Item(186, 1).Value2 = 9
$ws.Cells.Item(186, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(186, 3).Value2 = "Metropolitana"
$ws.Cells.Item(186, 4).NumberFormat = $dateFmt
$ws.Cells.Item(186, 4).Value2 = 45258
$ws.Cells.Item(186, 5).Value2 = 13
$ws.Cells.Item(186, 6).Value2 = 300000000
$ws.Cells.Item(186, 7).Value2 = "Espárragos"
$ws.Cells.Item(186, 8).Value2 = "Sin especificar"
$ws.Cells.Item(186, 9).Value2 = "Banquete"
$ws.Cells.Item(186, 10).Value2 = 300
$ws.Cells.Item(186, 11).Value2 = 1800
$ws.Cells.Item(186, 12).Value2 = 1800
$ws.Cells.Item(186, 13).Value2 = 1800
$ws.Cells.Item(186, 14).Value2 = "`$/kilo"
$ws.Cells.Item(186, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(186, 16).Value2 = 1800
$ws.Cells.Item(186, 17).Value2 = 1
$ws.Cells.Item(186, 18).Value2 = "Hortaliza"

# New row 187: Espárragos, Primera
$ws.Cells.Item(187, 1).Value2 = 9
$ws.Cells.Item(187, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(187, 3).Value2 = "Metropolitana"
$ws.Cells.Item(187, 4).NumberFormat = $dateFmt
$ws.Cells.Item(187, 4).Value2 = 45258
$ws.Cells.Item(187, 5).Value2 = 13
$ws.Cells.Item(187, 6).Value2 = 300000000
$ws.Cells.Item(187, 7).Value2 = "Espárragos"
$ws.Cells.Item(187, 8).Value2 = "Sin especificar"
$ws.Cells.Item(187, 9).Value2 = "Primera"
$ws.Cells.Item(187, 10).Value2 = 450
$ws.Cells.Item(187, 11).Value2 = 1600
$ws.Cells.Item(187, 12).Value2 = 1600
$ws.Cells.Item(187, 13).Value2 = 1600
$ws.Cells.Item(187, 14).Value2 = "`$/kilo"
$ws.Cells.Item(187, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(187, 16).Value2 = 1600
$ws.Cells.Item(187, 17).Value2 = 1
$ws.Cells.Item(187, 18).Value2 = "Hortaliza"

# New row 188: Espárragos, Segunda
$ws.Cells.Item(188, 1).Value2 = 9
$ws.Cells.Item(188, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(188, 3).Value2 = "Metropolitana"
$ws.Cells.Item(188, 4).NumberFormat = $dateFmt
$ws.Cells.Item(188, 4).Value2 = 45258
$ws.Cells.Item(188, 5).Value2 = 13
$ws.Cells.Item(188, 6).Value2 = 300000000
$ws.Cells.Item(188, 7).Value2 = "Espárragos"
$ws.Cells.Item(188, 8).Value2 = "Sin especificar"
$ws.Cells.Item(188, 9).Value2 = "Segunda"
$ws.Cells.Item(188, 10).Value2 = 200
$ws.Cells.Item(188, 11).Value2 = 1400
$ws.Cells.Item(188, 12).Value2 = 1400
$ws.Cells.Item(188, 13).Value2 = 1400
$ws.Cells.Item(188, 14).Value2 = "`$/kilo"
$ws.Cells.Item(188, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(188, 16).Value2 = 1400
$ws.Cells.Item(188, 17).Value2 = 1
$ws.Cells.Item(188, 18).Value2 = "Hortaliza"
